$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = 18
$v_B19 = @'
Português
'@
$ws.Range("B19").Value = $v_B19
$v_C19 = @'
Travessão
'@
$ws.Range("C19").Value = $v_C19
$v_D19 = @'
Classificação e significado da palavra <b>dúbio</b>
'@
$ws.Range("D19").Value = $v_D19
$v_E19 = @'
<b>adjetivo</b>
<ul>
	<li>sujeito a diferentes interpretações; ambíguo.</li>
	<li>difícil de caracterizar; impreciso, indefinível, vago.</li>
</ul>
'@
$ws.Range("E19").Value = $v_E19
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0

# Row 20
$ws.Range("A20").Value = 19
$v_B20 = @'
Português
'@
$ws.Range("B20").Value = $v_B20
$v_C20 = @'
Sinônimo E Antônimo
'@
$ws.Range("C20").Value = $v_C20
$v_D20 = @'
Classficação e significado da palavra <b>reminiscência</b>
'@
$ws.Range("D20").Value = $v_D20
$v_E20 = @'
<b>substantivo feminino</b>
<ol>
	<li>imagem lembrada do passado; o que se conserva na memória.</li>
	<li>lembrança vaga ou incompleta.</li>
</ol>
'@
$ws.Range("E20").Value = $v_E20
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0

# Row 21
$ws.Range("A21").Value = 20
$v_B21 = @'
Português
'@
$ws.Range("B21").Value = $v_B21
$v_C21 = @'
Sinônimo E Antônimo
'@
$ws.Range("C21").Value = $v_C21
$v_D21 = @'
Classificação e significado da palavra <b>resignação</b>
'@
$ws.Range("D21").Value = $v_D21
$v_E21 = @'
<b>substantivo feminino</b>
<ol>
	<li>submissão à vontade de alguém ou ao destino.</li>
	<li>demissão voluntária de um cargo.</li>
</ol>
'@
$ws.Range("E21").Value = $v_E21
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0

# Row 22
$ws.Range("A22").Value = 21
$v_B22 = @'
Inglês
'@
$ws.Range("B22").Value = $v_B22
$v_C22 = @'
Semantic
'@
$ws.Range("C22").Value = $v_C22
$v_D22 = @'
classificação e significado da palavra, <i>em inglês</i>, <b>Hence</b>
'@
$ws.Range("D22").Value = $v_D22
$v_E22 = @'
<b>Adverb</b>
Por isso
'@
$ws.Range("E22").Value = $v_E22
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0

# Row 23
$ws.Range("A23").Value = 22
$v_B23 = @'
Inglês
'@
$ws.Range("B23").Value = $v_B23
$v_C23 = @'
Semantic
'@
$ws.Range("C23").Value = $v_C23
$v_D23 = @'
Classificação e significado da palavra. <i>em inglês</i>, <b>yield</b>
'@
$ws.Range("D23").Value = $v_D23
$v_E23 = @'
<b>Nouns</b>
<ul>
	<li>rendimento</li>
	<li>produção</li>
</ul>
<b>Verb</b>
<ul>
	<li>produzir</li>
	<li>dar</li>
</ul>
'@
$ws.Range("E23").Value = $v_E23
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
